$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the used range with a new "Save" column (H), matching the
# formatting (bold/border/centered) already used by the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
